$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reframed_growth_rates")
$ws.Activate()

# --- Update data rows (A2:E12) ---
# Row 2: CDM35 Glucose / +
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.24
$ws.Range("E2").Value = 0.15

# Row 3: CDM35 Glucose / -
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.24
$ws.Range("E3").Value = 0.15

# Row 4: CMD35 - Glucose / +
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5: CDM35 Lactose / +
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.28999999999999998

# Row 6: CDM35 Lactose / "- " (note trailing space - new label variant)
$ws.Range("B6").Value = "- "
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.21

# Row 7: CDM35 Galactose / +
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.24
$ws.Range("E7").Value = 0.15

# Row 8: CDM35 Galactose / "- " (note trailing space - new label variant)
$ws.Range("B8").Value = "- "
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

# Row 9: CDM42 Glucose / +
$ws.Range("C9").Value = 1.18
$ws.Range("D9").Value = 0.25
$ws.Range("E9").Value = 0.28000000000000003

# Row 10: CDM42 - Glucose / +
$ws.Range("C10").Value = 0.92
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.01

# Row 11: CDM42  Galactose / +
$ws.Range("A11").Value = "CDM42  Galactose"
$ws.Range("C11").Value = 0.47
$ws.Range("D11").Value = 0.25
$ws.Range("E11").Value = 0.28000000000000003

# Row 12: CDM42  Lactose / +
$ws.Range("A12").Value = "CDM42  Lactose"
$ws.Range("C12").Value = 1.44
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.48

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 17.166666666666668
$ws.Columns.Item(2).ColumnWidth = 15.373697916666666

# --- Selection ---
$ws.Range("F19").Select() | Out-Null
